# Swap the B:G (Code, Rate, Rate, Qty, Value) data between pairs of rows.
# The underlying report re-sorted/re-paired certain stock-lot rows; the
# serial number (col A) and item description (col C) stay attached to
# their original row, while the lot code / rates / qty / value (B, D, E,
# F, G) move to the other row in the pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(313, 314),
    @(316, 318),
    @(351, 352),
    @(355, 356),
    @(372, 373),
    @(375, 376),
    @(389, 390),
    @(400, 401),
    @(419, 420),
    @(421, 422),
    @(579, 580),
    @(590, 591),
    @(593, 594),
    @(604, 605),
    @(687, 688),
    @(709, 710),
    @(720, 721),
    @(859, 860)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Columns B (2) through G (7): Code, Rate, Rate, Qty, Value
    for ($c = 2; $c -le 7; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

Write-Host "Swapped $($rowPairs.Count) row pairs"
